$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    ,@(2, 'Bitcoin', 'https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc', '27.549.03', '  -1.35%  ', $false)
    ,@(3, 'Ethereum', 'https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth', '1.843.38', '  -2.24%  ', $false)
    ,@(4, 'TetherUSD', 'https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt', '1.006', '  -1.16%  ', $true)
    ,@(5, 'BNB', 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb', '333.78', '  -0.58%  ', $true)
    ,@(6, 'USDC', 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc', '1.006', '  -1.08%  ', $true)
    ,@(7, 'XRP', 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp', '0.4628', '  -1.32%  ', $true)
    ,@(8, 'Cardano', 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada', '0.3853', '  -1.88%  ', $true)
    ,@(9, 'OKB', 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb', '45.86', '  -2.59%  ', $true)
    ,@(10, 'Dogecoin', 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge', '0.07908', '  -0.67%  ', $true)
    ,@(11, 'Polygon', 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic', '0.9948', '  -2.00%  ', $true)
    ,@(12, 'Solana', 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol', '21.48', '  -0.97%  ', $true)
    ,@(13, 'WrappedEther', 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth', '1.859.30', '  -1.86%  ', $false)
    ,@(14, 'Polkadot', 'https://coinranking.com/coin/25W7FG7om+polkadot-dot', '5.927', '  -0.60%  ', $true)
    ,@(15, 'Chainlink', 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link', '7.119', '  -0.12%  ', $true)
    ,@(16, 'BinanceUSD', 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd', '1.007', '  -1.15%  ', $true)
    ,@(17, 'Litecoin', 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc', '88.71', '  +1.38%  ', $true)
    ,@(18, 'TRON', 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx', '0.06682', '  -1.61%  ', $true)
    ,@(19, 'ShibaInu', 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib', '0.00001034', '  -1.21%  ', $true)
    ,@(20, 'Avalanche', 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax', '17.07', '  +0.25%  ', $true)
    ,@(21, 'Dai', 'https://coinranking.com/coin/MoTuySvg7+dai-dai', '1.006', '  -1.12%  ', $true)
    ,@(22, 'WrappedBTC', 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc', '27.555.18', '  -1.35%  ', $false)
    ,@(23, 'Uniswap', 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni', '5.380', '  -1.82%  ', $true)
    ,@(24, 'Cosmos', 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom', '10.89', '  -0.86%  ', $true)
    ,@(25, 'Toncoin', 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton', '2.308', '  -2.03%  ', $true)
    ,@(26, 'Monero', 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr', '158.78', '  -0.62%  ', $true)
    ,@(27, 'EthereumClassic', 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc', '19.48', '  -2.76%  ', $true)
    ,@(28, 'LidoDAOToken', 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo', '2.103', '  +1.08%  ', $true)
    ,@(29, 'InternetComputer(DFINITY)', 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp', '5.398', '  -1.29%  ', $true)
    ,@(30, 'BitcoinCash', 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch', '119.94', '  -0.95%  ', $true)
    ,@(31, 'ImmutableX', 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx', '0.9744', '  +1.38%  ', $true)
    ,@(32, 'Stellar', 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm', '0.09381', '  -1.80%  ', $true)
    ,@(33, 'HuobiToken', 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht', '3.594', '  -1.73%  ', $true)
    ,@(34, 'Filecoin', 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil', '5.272', '  -1.10%  ', $true)
    ,@(35, 'ARBITRUM', 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb', '1.338', '  -1.20%  ', $true)
    ,@(36, 'Hedera', 'https://coinranking.com/coin/jad286TjB+hedera-hbar', '0.06027', '  -1.74%  ', $true)
    ,@(37, 'VeChain', 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet', '0.02226', '  -0.61%  ', $true)
    ,@(38, 'FraxShare', 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs', '8.284', '  +1.65%  ', $true)
    ,@(39, 'TrustWalletToken', 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt', '1.183', '  -1.91%  ', $true)
    ,@(40, 'TheSandbox', 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand', '0.5878', '  -0.20%  ', $true)
    ,@(41, 'Algorand', 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo', '0.1865', '  -1.73%  ', $true)
    ,@(42, 'Aptos', 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt', '10.30', '  +0.88%  ', $true)
    ,@(43, 'WEMIXTOKEN', 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix', '1.239', '  -2.51%  ', $true)
    ,@(44, 'Decentraland', 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana', '0.5577', '  -1.18%  ', $true)
    ,@(45, 'EnergySwap', 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens', '12.18', '  +0.57%  ', $true)
    ,@(46, 'NEARProtocol', 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near', '1.902', '  -1.22%  ', $true)
    ,@(47, 'Cronos', 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro', '0.06700', '  -2.26%  ', $true)
    ,@(48, 'Quant', 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt', '111.04', '  -2.24%  ', $true)
    ,@(49, 'EOS', 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos', '1.051', '  -1.46%  ', $true)
    ,@(50, 'PaxDollar', 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp', '1.006', '  -1.22%  ', $true)
    ,@(51, 'Aave', 'https://coinranking.com/coin/ixgUfzmLR+aave-aave', '70.09', '  -1.05%  ', $true)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Range("B$r").Value = $row[1]
    $ws.Range("C$r").Value = $row[2]

    $dCell = $ws.Range("D$r")
    if ($row[5]) {
        # Numeric-looking price strings must stay text; force with a quote prefix
        # and then restore the cell's style so no stray number format lingers.
        $dCell.Value = "'" + $row[3]
        $dCell.Style = "Normal"
    } else {
        $dCell.Value = $row[3]
    }

    $ws.Range("E$r").Value = $row[4]
}
